# Adds a "Run 50" column of data, pushing the existing "Mean" column
# one position to the right (from AZ to BA) and filling the vacated
# AZ column with the new "Run 50" values.
#
# Column AZ (52) currently holds the "Mean" header/values; it becomes the
# "Run 50" header/values. A brand-new column BA (53) is appended to hold
# the (relocated) "Mean" header/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy AZ1's formatting (bold/border/center header style) into the new
# BA1 cell before overwriting the text, so BA1 ends up styled exactly
# like every other header cell.
$ws.Range("AZ1").Copy($ws.Range("BA1"))
$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

# New data values, identical for every data row (2 through 14), matching
# the pattern already used by every other run/mean column in this sheet.
$newRunValue = 604.24170852
$newMeanValue = 852.41534517

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 52).Value = $newRunValue
    $ws.Cells.Item($row, 53).Value = $newMeanValue
}
